$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item("TableForCopy")

$ws.Range("M18").FormulaArray = "=NameError"
$lo.ShowTotals = $true
$ws.Range("N18").Formula = "=NA()"

$lo.ShowTotals = $false
$ws.Range("M19").Value = "'=foo"
$ws.Range("N19").Value = 100

$lo.ShowTotals = $true
Write-Host "ref after re-enable totals:" $lo.Range.Address()
